# Add last month's data in details table: update the percentage
# distributions in columns C-H (the underlying column totals still sum
# to 100 after the update, reflecting the refreshed monthly figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3415863323499416
$ws.Range("H2").Value = 0.3512440547385684

$ws.Range("G3").Value = 1.980210622318502
$ws.Range("H3").Value = 1.962257289042282

$ws.Range("G4").Value = 3.749079902218471
$ws.Range("H4").Value = 2.838173106021388

$ws.Range("G5").Value = 4.072904833662385
$ws.Range("H5").Value = 2.289452872001473

$ws.Range("G6").Value = 11.84997093778064
$ws.Range("H6").Value = 13.0843329994615

$ws.Range("G7").Value = 27.67922478487725
$ws.Range("H7").Value = 29.43473191536076

$ws.Range("G8").Value = 50.32702258679282
$ws.Range("H8").Value = 50.03980776337402

$ws.Range("E9").Value = 0.03650486079923238
$ws.Range("F9").Value = 0.09072711310134488

$ws.Range("E10").Value = 0.0989304403502218
$ws.Range("F10").Value = 0.1728394345701645

$ws.Range("E11").Value = 0.2506919250226087
$ws.Range("F11").Value = 0.2734730094815129

$ws.Range("E12").Value = 1.559117991228031
$ws.Range("F12").Value = 1.635404270469866

$ws.Range("E13").Value = 32.72254985628717
$ws.Range("F13").Value = 32.16493109235962

$ws.Range("E14").Value = 65.33220492631274
$ws.Range("F14").Value = 65.66262508001749

$ws.Range("C15").Value = 0.397506286674684
$ws.Range("D15").Value = 0.5719799423698944

$ws.Range("C16").Value = 1.445152178463545
$ws.Range("D16").Value = 1.470712146025029

$ws.Range("C17").Value = 2.169802081809555
$ws.Range("D17").Value = 2.116704947688595

$ws.Range("C18").Value = 3.955652293326212
$ws.Range("D18").Value = 3.361627762342925

$ws.Range("C19").Value = 6.682533827737573
$ws.Range("D19").Value = 7.727861128481588

$ws.Range("C20").Value = 19.57661370401237
$ws.Range("D20").Value = 21.25909838837807

$ws.Range("C21").Value = 65.77273962797607
$ws.Range("D21").Value = 63.4920156847139
